# Apply the edit described in the diff:
#  - Insert two new data rows right before the current row 467, shifting all
#    existing data rows (467-564) down by two rows (to 469-566).
#  - Populate the two newly inserted rows (467 and 468) with new data.
#  - The worksheet dimension will automatically grow from A1:T564 to A1:T566.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 467, pushing everything at/after row 467 down by 2.
$ws.Rows.Item(467).Resize(2).Insert()

# New row 467: "1a amarillo" entry (Provincia de Melipilla, $/malla 18 kilos)
$ws.Range("A467").Value = 4
$ws.Range("B467").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C467").Value = "Los Lagos"
$ws.Range("D467").Value = 44782
$ws.Range("E467").Value = 10
$ws.Range("F467").Value = "Fruta"
$ws.Range("G467").Value = 100102
$ws.Range("H467").Value = "Cítricos"
$ws.Range("I467").Value = 100102003
$ws.Range("J467").Value = "Limón"
$ws.Range("K467").Value = "Sin especificar"
$ws.Range("L467").Value = "1a amarillo"
$ws.Range("M467").Value = 1400
$ws.Range("N467").Value = 8500
$ws.Range("O467").Value = 9000
$ws.Range("P467").Value = 8750
$ws.Range("Q467").Value = "$/malla 18 kilos"
$ws.Range("R467").Value = "Provincia de Melipilla"
$ws.Range("S467").Value = 486
$ws.Range("T467").Value = 18

# New row 468: "2a amarillo" entry (Provincia de Melipilla, $/malla 18 kilos)
$ws.Range("A468").Value = 4
$ws.Range("B468").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C468").Value = "Los Lagos"
$ws.Range("D468").Value = 44782
$ws.Range("E468").Value = 10
$ws.Range("F468").Value = "Fruta"
$ws.Range("G468").Value = 100102
$ws.Range("H468").Value = "Cítricos"
$ws.Range("I468").Value = 100102003
$ws.Range("J468").Value = "Limón"
$ws.Range("K468").Value = "Sin especificar"
$ws.Range("L468").Value = "2a amarillo"
$ws.Range("M468").Value = 700
$ws.Range("N468").Value = 6500
$ws.Range("O468").Value = 6500
$ws.Range("P468").Value = 6500
$ws.Range("Q468").Value = "$/malla 18 kilos"
$ws.Range("R468").Value = "Provincia de Melipilla"
$ws.Range("S468").Value = 361
$ws.Range("T468").Value = 18
